# Generate Report for Handoff
# - flips "In Translation" rows to "Ready for handoff"
# - refreshes the HO Xliff generate / handoff timestamps
# - widens the status-adjacent date columns to fit the new values

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status: "In Translation" -> "Ready for handoff" ---
$wsOverview.Cells.Item(2, 5).Value = "Ready for handoff"   # E2 (zh-cn column)
$wsOverview.Cells.Item(2, 6).Value = "Ready for handoff"   # F2 (de-de column)
$wsZhCn.Cells.Item(2, 3).Value = "Ready for handoff"       # C2 Status
$wsDeDe.Cells.Item(2, 3).Value = "Ready for handoff"       # C2 Status

# --- Timestamps ---
$wsOverview.Cells.Item(2, 7).Value = "2016-09-07 03:14:51" # G2 Latest HO Xliff Generate Date
$wsOverview.Cells.Item(2, 7).NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsDeDe.Cells.Item(2, 8).Value = "2016-09-07 03:14:51"     # H2 Latest Handoff Datetime
$wsDeDe.Cells.Item(2, 8).NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsZhCn.Cells.Item(2, 8).Value = "2016-09-07 03:14:46"     # H2 Latest Handoff Datetime
$wsZhCn.Cells.Item(2, 8).NumberFormat = "yyyy-mm-dd HH:mm:ss"

# --- Column widths (status/date columns got a bit wider) ---
$wsOverview.Columns.Item(5).ColumnWidth = 16.3333333333333
$wsOverview.Columns.Item(6).ColumnWidth = 16.3333333333333
$wsZhCn.Columns.Item(3).ColumnWidth = 16.3333333333333
$wsDeDe.Columns.Item(3).ColumnWidth = 16.3333333333333
